# Generate Report for Handoff
# Adds a new tracked file (aff6cbec-58b3-460e-b055-9fc8d82fe812.md) as row 9
# to the Overview / zh-cn / de-de sheets of the localization-status workbook.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Expand the three tables by one row so refs/dimension/autofilter grow ----
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$fileBase   = "aff6cbec-58b3-460e-b055-9fc8d82fe812"
$fileName   = "$fileBase.md"
$pathName   = "e2e\$fileBase.md"
$zhXlfName  = "$fileBase.095c20b14c39c7a7cdcffbe903e0ced8b1560904.zh-cn.xlf"
$deXlfName  = "$fileBase.095c20b14c39c7a7cdcffbe903e0ced8b1560904.de-de.xlf"
$handoffDt  = "2016-08-19 02:41:33"
$zhHandoffDt= "2016-08-19 02:41:29"

# --------------------------- Overview sheet row 9 ---------------------------
$wsOverview.Range("A9").Value = $fileName
$wsOverview.Range("B9").Value = $pathName
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = "'"
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = $handoffDt
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$fileBase/e2e/$fileName",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $pathName
) | Out-Null
$wsOverview.Range("B9").Style = "HyperLink"

# ----------------------------- zh-cn sheet row 9 -----------------------------
$wsZhCn.Range("A9").Value = $fileName
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "'False"
$wsZhCn.Range("G9").Value = $zhXlfName
$wsZhCn.Range("H9").Value = $zhHandoffDt
$wsZhCn.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I9").Value = "'"
$wsZhCn.Range("J9").Value = "'"
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L9").Value = "'"
$wsZhCn.Range("M9").Value = "'True"
$wsZhCn.Range("N9").Value = "'"
$wsZhCn.Range("O9").Value = "'False"
$wsZhCn.Range("P9").Value = "'"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$fileBase/e2e/$fileName",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $fileName
) | Out-Null
$wsZhCn.Range("A9").Style = "HyperLink"

# ----------------------------- de-de sheet row 9 -----------------------------
$wsDeDe.Range("A9").Value = $fileName
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "'False"
$wsDeDe.Range("G9").Value = $deXlfName
$wsDeDe.Range("H9").Value = $handoffDt
$wsDeDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I9").Value = "'"
$wsDeDe.Range("J9").Value = "'"
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L9").Value = "'"
$wsDeDe.Range("M9").Value = "'True"
$wsDeDe.Range("N9").Value = "'"
$wsDeDe.Range("O9").Value = "'False"
$wsDeDe.Range("P9").Value = "'"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$fileBase/e2e/$fileName",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $fileName
) | Out-Null
$wsDeDe.Range("A9").Style = "HyperLink"
